$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both contain the same table with column F
# ("想去人数") values that need updating.
$sheetNames = @("展览", "全部类型")

$updates = @{
    "F3" = 2412
    "F4" = 450
    "F5" = 84
    "F6" = 6505
    "F7" = 351
    "F8" = 130
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    foreach ($cell in $updates.Keys) {
        $ws.Range($cell).Value = $updates[$cell]
    }
}
